$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.185199999999993
$ws.Range("D21").Value = -8.619099999999998
$ws.Range("D23").Value = -7.218999999999997
$ws.Range("D25").Value = -8.988899999999992
$ws.Range("E27").Value = 16.57029999999999
$ws.Range("E31").Value = 16.2095
$ws.Range("E39").Value = 15.6731
$ws.Range("E48").Value = 17.2533
$ws.Range("E51").Value = 17.1219
$ws.Range("E52").Value = 16.9359
$ws.Range("D53").Value = -6.140499999999998
$ws.Range("E55").Value = 16.6945
$ws.Range("E56").Value = 16.02740000000001
$ws.Range("D57").Value = -8.096699999999991
$ws.Range("E57").Value = 16.547
$ws.Range("D59").Value = -8.440099999999999
$ws.Range("D69").Value = -7.101999999999997
$ws.Range("E73").Value = 17.4558
$ws.Range("D79").Value = -6.025999999999997
$ws.Range("D83").Value = -8.910199999999996
$ws.Range("E89").Value = 17.19380000000002
$ws.Range("E90").Value = 16.20139999999999
$ws.Range("D93").Value = -6.355999999999995
